{"js": "// Office.js (Word JavaScript API) script.\n// Applies the \"Added many more features\" edits to the Football Star review:\n//  - Updates the title (appears twice: H1 heading + bold recap line)\n//  - Rewrites the four \"What we like\" bullets\n//  - Rewrites the two \"What we don't like\" bullets\n//  - Rewrites the italic meta-description line\n//\n// Implemented as straight text replacements via Body.search so the existing\n// run formatting (bold/italic/heading styles) is preserved untouched.\n\nconst replacements = [\n  [\n    \"Play Football Star for Free: Exciting Soccer-Themed Slot Game\",\n    \"Play Football Star Free - A Thrilling Soccer-Themed Online Slot Game\",\n  ],\n  [\n    \"Fun soccer-themed design and vibrant graphics\",\n    \"Straightforward gameplay mechanics\",\n  ],\n  [\n    \"Several special features, including free spins and striking wilds\",\n    \"Special features that enhance the gameplay\",\n  ],\n  [\n    \"Rolling reels feature increases chances of winning\",\n    \"Visually appealing soccer-themed design\",\n  ],\n  [\n    \"Developed by reputable game developer, Microgaming\",\n    \"Developed by reputable online casino game developer\",\n  ],\n  [\n    \"No progressive jackpot feature\",\n    \"Limited number of paylines\",\n  ],\n  [\n    \"No gamble feature to increase winnings\",\n    \"May not appeal to non-soccer fans\",\n  ],\n  [\n    \"Review of Football Star, an exciting soccer-themed online slot game from Microgaming. Play for free with rolling reels, stacked wilds, and a free spins bonus round.\",\n    \"Play Football Star free and experience the excitement of soccer on the reels.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the \"Added many more features\" edits to the Football Star review:\n#  - Updates the title (appears twice: H1 heading + bold recap line)\n#  - Rewrites the four \"What we like\" bullets\n#  - Rewrites the two \"What we don't like\" bullets\n#  - Rewrites the italic meta-description line\n#\n# Implemented with Find/Replace (wdReplaceAll) over the whole document\n# range so existing run formatting (bold/italic/heading styles) is left\n# untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Football Star for Free: Exciting Soccer-Themed Slot Game\", \"Play Football Star Free - A Thrilling Soccer-Themed Online Slot Game\"),\n    @(\"Fun soccer-themed design and vibrant graphics\", \"Straightforward gameplay mechanics\"),\n    @(\"Several special features, including free spins and striking wilds\", \"Special features that enhance the gameplay\"),\n    @(\"Rolling reels feature increases chances of winning\", \"Visually appealing soccer-themed design\"),\n    @(\"Developed by reputable game developer, Microgaming\", \"Developed by reputable online casino game developer\"),\n    @(\"No progressive jackpot feature\", \"Limited number of paylines\"),\n    @(\"No gamble feature to increase winnings\", \"May not appeal to non-soccer fans\"),\n    @(\"Review of Football Star, an exciting soccer-themed online slot game from Microgaming. Play for free with rolling reels, stacked wilds, and a free spins bonus round.\", \"Play Football Star free and experience the excitement of soccer on the reels.\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
